$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.245.92'
$ws.Range('E2').Value = '  -0.28%  '

$ws.Range('D3').Value = '1.841.35'
$ws.Range('E3').Value = '  +0.00%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.89'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6693'
$ws.Range('E6').Value = '  -2.56%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.0000'
$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07428'
$ws.Range('E8').Value = '  -1.39%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2961'
$ws.Range('E9').Value = '  -2.35%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.85'
$ws.Range('E10').Value = '  -1.92%  '

$ws.Range('E11').Value = '  +0.73%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.033'
$ws.Range('E12').Value = '  -1.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6792'
$ws.Range('E13').Value = '  -1.21%  '

$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '86.34'
$ws.Range('E14').Value = '  -3.73%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.678.21'
$ws.Range('E15').Value = '  -8.86%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.200'
$ws.Range('E16').Value = '  -1.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008269'
$ws.Range('E17').Value = '  +0.08%  '

$ws.Range('D18').Value = '28.917.77'
$ws.Range('E18').Value = '  -1.34%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.96'
$ws.Range('E19').Value = '  -2.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.55'
$ws.Range('E20').Value = '  -0.38%  '

$ws.Range('E21').Value = '  -0.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.247'
$ws.Range('E22').Value = '  -3.27%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.23'
$ws.Range('E24').Value = '  +0.08%  '

$ws.Range('E25').Value = '  -1.60%  '

$ws.Range('E26').Value = '  -2.86%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.02'
$ws.Range('E27').Value = '  -0.62%  '

$ws.Range('E28').Value = '  -1.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.204'
$ws.Range('E29').Value = '  -0.72%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.086'
$ws.Range('E30').Value = '  -1.31%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.200'
$ws.Range('E31').Value = '  -0.13%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05360'
$ws.Range('E32').Value = '  +4.26%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7580'
$ws.Range('E33').Value = '  -1.78%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.873'
$ws.Range('E34').Value = '  +1.18%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.140'
$ws.Range('E35').Value = '  -0.06%  '

$ws.Range('D37').Value = '1.331.73'
$ws.Range('E37').Value = '  +2.67%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01806'
$ws.Range('E38').Value = '  -2.13%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.731'
$ws.Range('E39').Value = '  +1.37%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9210'
$ws.Range('E40').Value = '  -2.46%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.005'
$ws.Range('E41').Value = '  +6.59%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  +0.11%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.46'
$ws.Range('E43').Value = '  -2.09%  '

$ws.Range('E44').Value = '  +4.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07799'
$ws.Range('E45').Value = '  +12.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5163'
$ws.Range('E46').Value = '  -0.84%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.767'
$ws.Range('E47').Value = '  -0.17%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.78'
$ws.Range('E48').Value = '  +0.76%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '1.915.94'
$ws.Range('E49').Value = '  -3.72%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.280'
$ws.Range('E50').Value = '  -4.24%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05939'
$ws.Range('E51').Value = '  +0.23%  '
